$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-12 18:56:33"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-12 18:56:26"
$wsZhCn.Range("K4").Value = "2016-08-12 18:56:55"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-12 18:56:33"
$wsDeDe.Range("K4").Value = "2016-08-12 18:57:11"
